$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2234.52
$ws.Range("J17").Value = 2234.52
$ws.Range("L17").Value = 6703.559999999999
$ws.Range("N17").Value = -7039.559999999999
$ws.Range("H74").Value = 7767.8423
$ws.Range("I74").Value = 4570
$ws.Range("J74").Value = 9633.25
$ws.Range("K74").Value = 4570
$ws.Range("L74").Value = 9633.25
$ws.Range("M74").Value = -3634
$ws.Range("N74").Value = -11505.25
$ws.Range("H77").Value = 7767.8423
$ws.Range("I77").Value = 4570
$ws.Range("J77").Value = 9633.25
$ws.Range("K77").Value = 22850
$ws.Range("L77").Value = 48166.25
$ws.Range("M77").Value = -18170
$ws.Range("N77").Value = -57526.25
$ws.Range("H86").Value = 2812.125
$ws.Range("I86").Value = 1013.46155
$ws.Range("K86").Value = 1013.46155
$ws.Range("M86").Value = 109.53845
$ws.Range("H89").Value = 2812.125
$ws.Range("I89").Value = 1013.46155
$ws.Range("K89").Value = 5067.30775
$ws.Range("M89").Value = 548.6922500000001
$ws.Range("H112").Value = 1463.4889
$ws.Range("J112").Value = 1491.7906
$ws.Range("L112").Value = 4475.3718
$ws.Range("N112").Value = -6691.3718
$ws.Range("H132").Value = 3243.818
$ws.Range("I132").Value = 3243.818
$ws.Range("K132").Value = 9731.454000000002
$ws.Range("M132").Value = -7201.454000000002
$ws.Range("H138").Value = 3472.5186
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3472.5186
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10417.5558
$ws.Range("N138").Value = -20697.5558
$ws.Range("M138").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8277.727999999999
$ws.Range("I32").Value = 7221.6313
$ws.Range("K32").Value = 7221.6313
$ws.Range("M32").Value = -6934.6313
$ws.Range("H61").Value = 7187.303
$ws.Range("I61").Value = 7161.476
$ws.Range("J61").Value = 7232.5
$ws.Range("K61").Value = 7161.476
$ws.Range("L61").Value = 7232.5
$ws.Range("M61").Value = -6949.476
$ws.Range("N61").Value = -7656.5
$ws.Range("H132").Value = 4774.8486
$ws.Range("J132").Value = 13999
$ws.Range("L132").Value = 41997
$ws.Range("N132").Value = -47057
$ws.Range("H136").Value = 7187.303
$ws.Range("I136").Value = 7161.476
$ws.Range("J136").Value = 7232.5
$ws.Range("K136").Value = 21484.428
$ws.Range("L136").Value = 21697.5
$ws.Range("M136").Value = -18934.428
$ws.Range("N136").Value = -26797.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1898.2
$ws.Range("I20").Value = 1899.3636
$ws.Range("J20").Value = 1896.7778
$ws.Range("K20").Value = 1899.3636
$ws.Range("L20").Value = 1896.7778
$ws.Range("M20").Value = -1652.3636
$ws.Range("N20").Value = -2390.7778
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
$ws.Range("H88").Value = 27001.154
$ws.Range("J88").Value = 27001.154
$ws.Range("L88").Value = 27001.154
$ws.Range("N88").Value = -27813.154
$ws.Range("H91").Value = 27001.154
$ws.Range("J91").Value = 27001.154
$ws.Range("L91").Value = 27001.154
$ws.Range("N91").Value = -29809.154
$ws.Range("H103").Value = 24999.5
$ws.Range("J103").Value = 24999.5
$ws.Range("L103").Value = 24999.5
$ws.Range("N103").Value = -27343.5
$ws.Range("H115").Value = 231420.5
$ws.Range("J115").Value = 231420.5
$ws.Range("L115").Value = 231420.5
$ws.Range("N115").Value = -234554.5
$ws.Range("H116").Value = 53749.5
$ws.Range("J116").Value = 53749.5
$ws.Range("L116").Value = 53749.5
$ws.Range("N116").Value = -62927.5
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5524.246
$ws.Range("I31").Value = 2883.7693
$ws.Range("J31").Value = 7485.7427
$ws.Range("K31").Value = 2883.7693
$ws.Range("L31").Value = 7485.7427
$ws.Range("M31").Value = -2588.7693
$ws.Range("N31").Value = -8075.7427
$ws.Range("H34").Value = 5524.246
$ws.Range("I34").Value = 2883.7693
$ws.Range("J34").Value = 7485.7427
$ws.Range("K34").Value = 2883.7693
$ws.Range("L34").Value = 7485.7427
$ws.Range("M34").Value = -2681.7693
$ws.Range("N34").Value = -7889.7427
$ws.Range("H62").Value = 8498.75
$ws.Range("J62").Value = 9665.333000000001
$ws.Range("L62").Value = 9665.333000000001
$ws.Range("N62").Value = -10913.333
$ws.Range("H65").Value = 8498.75
$ws.Range("J65").Value = 9665.333000000001
$ws.Range("L65").Value = 48326.665
$ws.Range("N65").Value = -54566.665
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H132").Value = 12604.9
$ws.Range("I132").Value = 14507.375
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 43522.125
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -40992.125
$ws.Range("N132").Value = -20045
$ws.Range("H134").Value = 2068.3333
$ws.Range("I134").Value = 2019.1666
$ws.Range("J134").Value = 2166.6667
$ws.Range("K134").Value = 6057.4998
$ws.Range("L134").Value = 6500.000100000001
$ws.Range("M134").Value = -3522.4998
$ws.Range("N134").Value = -11570.0001
$ws.Range("N115").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 184182080
$ws.Range("J4").Value = 8333530.5
$ws.Range("K4").Value = 552546240
$ws.Range("L4").Value = 25000591.5
$ws.Range("M4").Value = -552546128
$ws.Range("N4").Value = -25000815.5
$ws.Range("H21").Value = 1083.3334
$ws.Range("J21").Value = 1450
$ws.Range("L21").Value = 4350
$ws.Range("N21").Value = -4696
$ws.Range("H76").Value = 7747.5
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 9330
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 27990
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -28756
$ws.Range("H79").Value = 7747.5
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 9330
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 27990
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -30642
$ws.Range("H80").Value = 4597.4443
$ws.Range("I80").Value = 3608.5
$ws.Range("K80").Value = 10825.5
$ws.Range("M80").Value = -9889.5
$ws.Range("H83").Value = 4597.4443
$ws.Range("I83").Value = 3608.5
$ws.Range("K83").Value = 32476.5
$ws.Range("M83").Value = -27796.5
$ws.Range("H131").Value = 1585.5
$ws.Range("J131").Value = 1897.2941
$ws.Range("L131").Value = 5691.8823
$ws.Range("N131").Value = -15771.8823

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 20000
$ws.Range("I38").Value = 20000
$ws.Range("K38").Value = 20000
$ws.Range("M38").Value = -19537
$ws.Range("H42").Value = 79998.5
$ws.Range("J42").Value = 79998.5
$ws.Range("L42").Value = 79998.5
$ws.Range("N42").Value = -80968.5
$ws.Range("H115").Value = 79998.5
$ws.Range("J115").Value = 79998.5
$ws.Range("L115").Value = 79998.5
$ws.Range("N115").Value = -82348.5
$ws.Range("H132").Value = 8299.931
$ws.Range("I132").Value = 5790.7
$ws.Range("J132").Value = 13876
$ws.Range("K132").Value = 17372.1
$ws.Range("L132").Value = 41628
$ws.Range("M132").Value = -14842.1
$ws.Range("N132").Value = -46688

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 169526.83
$ws.Range("I22").Value = 1000000
$ws.Range("J22").Value = 3432.2
$ws.Range("K22").Value = 1000000
$ws.Range("L22").Value = 3432.2
$ws.Range("M22").Value = -999705
$ws.Range("N22").Value = -4022.2
$ws.Range("H27").Value = 169526.83
$ws.Range("I27").Value = 1000000
$ws.Range("J27").Value = 3432.2
$ws.Range("K27").Value = 1000000
$ws.Range("L27").Value = 3432.2
$ws.Range("M27").Value = -999893
$ws.Range("N27").Value = -3646.2
$ws.Range("H93").Value = 1846.0625
$ws.Range("I93").Value = 1535.875
$ws.Range("J93").Value = 2156.25
$ws.Range("K93").Value = 1535.875
$ws.Range("L93").Value = 2156.25
$ws.Range("M93").Value = -287.875
$ws.Range("N93").Value = -4652.25
